$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C9 and F9 currently hold numeric 0 formatted as a percentage (0%).
# The edit replaces them with the literal text "0%" (quote-prefixed text),
# matching how the author retyped the values as text in Excel - switching
# the cell format to General and entering the value with a leading
# apostrophe so it is stored as text rather than a number.
$ws.Range("C9").NumberFormat = "general"
$ws.Range("C9").Value = "'0%"

$ws.Range("F9").NumberFormat = "general"
$ws.Range("F9").Value = "'0%"

# Update the active selection to match the author's final cursor position.
$ws.Range("F9").Select()
